# Update Name of Algo
# Applies updated numeric results to Sheet1 of the workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D3").Value  = -7.811999999999999
$ws.Range("D4").Value  = -8.032
$ws.Range("E6").Value  = 12.718
$ws.Range("D7").Value  = -8.102
$ws.Range("E7").Value  = 12.895
$ws.Range("D8").Value  = -8.010999999999999
$ws.Range("E8").Value  = 13.03
$ws.Range("B11").Value = 6.289
$ws.Range("B12").Value = 6.066
$ws.Range("D12").Value = -8.442
$ws.Range("D14").Value = -8.224
$ws.Range("B15").Value = 6.336
$ws.Range("E19").Value = 12.169
$ws.Range("E21").Value = 13.136
$ws.Range("D22").Value = -7.812
$ws.Range("E24").Value = 12.718
$ws.Range("E25").Value = 12.169
